# Update cryptos list worksheet with latest scraped values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text columns (coin name, link, % change) - safe to assign directly
$ws.Range('E2').Value = '  +2.93%  '
$ws.Range('E3').Value = '  +3.51%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('E5').Value = '  +3.91%  '
$ws.Range('E6').Value = '  +2.62%  '
$ws.Range('E8').Value = '  +4.70%  '
$ws.Range('E9').Value = '  +3.31%  '
$ws.Range('E10').Value = '  +7.14%  '
$ws.Range('E11').Value = '  +3.22%  '
$ws.Range('E12').Value = '  +2.55%  '
$ws.Range('E13').Value = '  +3.54%  '
$ws.Range('E14').Value = '  +4.04%  '
$ws.Range('E15').Value = '  +6.79%  '
$ws.Range('E16').Value = '  +3.04%  '
$ws.Range('E17').Value = '  +3.09%  '
$ws.Range('E18').Value = '  -0.13%  '
$ws.Range('E19').Value = '  +4.56%  '
$ws.Range('E20').Value = '  +6.28%  '
$ws.Range('E21').Value = '  +4.33%  '
$ws.Range('E22').Value = '  +5.51%  '
$ws.Range('E23').Value = '  -0.18%  '
$ws.Range('E24').Value = '  -1.02%  '
$ws.Range('E25').Value = '  +7.82%  '
$ws.Range('E26').Value = '  +2.20%  '
$ws.Range('E27').Value = '  +5.28%  '
$ws.Range('E28').Value = '  +3.87%  '
$ws.Range('E29').Value = '  +4.15%  '
$ws.Range('E30').Value = '  +3.31%  '
$ws.Range('E31').Value = '  +7.80%  '
$ws.Range('E32').Value = '  +3.88%  '
$ws.Range('E33').Value = '  +6.34%  '
$ws.Range('E34').Value = '  +4.80%  '
$ws.Range('E35').Value = '  +2.50%  '
$ws.Range('E36').Value = '  +1.68%  '
$ws.Range('E37').Value = '  -0.96%  '
$ws.Range('E38').Value = '  +3.89%  '
$ws.Range('E39').Value = '  +6.55%  '
$ws.Range('E40').Value = '  +5.85%  '
$ws.Range('E41').Value = '  +3.18%  '
$ws.Range('E42').Value = '  -0.03%  '
$ws.Range('E44').Value = '  +3.17%  '
$ws.Range('B45').Value = 'BabyDogeCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('E45').Value = '  +3.63%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('E46').Value = '  +6.06%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('E47').Value = '  +3.00%  '
$ws.Range('B48').Value = 'Frax'
$ws.Range('C48').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('E48').Value = '  +0.18%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('E49').Value = '  +1.00%  '
$ws.Range('B50').Value = 'Aptos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('E50').Value = '  +5.95%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('E51').Value = '  +0.33%  '

# Price column - force text storage so formats like trailing zeros and
# "thousand.thousand.decimal" strings are preserved exactly instead of
# being auto-converted into floating point numbers by Excel.
$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '26.656.07'
$c.Style = 'Normal'
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '1.689.88'
$c.Style = 'Normal'
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '217.80'
$c.Style = 'Normal'
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '0.5343'
$c.Style = 'Normal'
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.2685'
$c.Style = 'Normal'
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '21.68'
$c.Style = 'Normal'
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.07800'
$c.Style = 'Normal'
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '1.679.50'
$c.Style = 'Normal'
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '0.5628'
$c.Style = 'Normal'
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '0.0₅8452'
$c.Style = 'Normal'
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '66.38'
$c.Style = 'Normal'
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '26.701.23'
$c.Style = 'Normal'
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '4.818'
$c.Style = 'Normal'
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '195.85'
$c.Style = 'Normal'
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '10.43'
$c.Style = 'Normal'
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '6.387'
$c.Style = 'Normal'
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '144.13'
$c.Style = 'Normal'
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '0.1294'
$c.Style = 'Normal'
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '7.493'
$c.Style = 'Normal'
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '16.29'
$c.Style = 'Normal'
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '1.426'
$c.Style = 'Normal'
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '0.06175'
$c.Style = 'Normal'
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '1.283'
$c.Style = 'Normal'
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '3.474'
$c.Style = 'Normal'
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '1.705'
$c.Style = 'Normal'
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '2.801'
$c.Style = 'Normal'
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '0.5743'
$c.Style = 'Normal'
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.01652'
$c.Style = 'Normal'
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '6.021'
$c.Style = 'Normal'
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '1.080.91'
$c.Style = 'Normal'
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.8653'
$c.Style = 'Normal'
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '1.840.55'
$c.Style = 'Normal'
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.0₈110'
$c.Style = 'Normal'
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '57.52'
$c.Style = 'Normal'
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '8.205'
$c.Style = 'Normal'
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '1.002'
$c.Style = 'Normal'
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '0.05224'
$c.Style = 'Normal'
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '6.114'
$c.Style = 'Normal'
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '0.4241'
$c.Style = 'Normal'
